# Insert a new weekly price record for "Macroferia Regional de Talca - Haba".
# This shifts the existing rows 53-62 down to 54-63 (preserving all of their
# data untouched) and populates the newly inserted row 53 with the new
# week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 53; Excel shifts rows 53:62 down to 54:63.
$ws.Rows.Item(53).Insert()

# Populate the newly inserted row 53 with the new record.
$ws.Range("A53").Value = 5
$ws.Range("B53").Value = "Macroferia Regional de Talca"
$ws.Range("C53").Value = "Maule"
$ws.Range("D53").Value = 44522
$ws.Range("E53").Value = 7
$ws.Range("F53").Value = 100112026
$ws.Range("G53").Value = "Haba"
$ws.Range("H53").Value = "Sin especificar"
$ws.Range("I53").Value = "Primera"
$ws.Range("J53").Value = 500
$ws.Range("K53").Value = 6000
$ws.Range("L53").Value = 6000
$ws.Range("M53").Value = 6000
$ws.Range("N53").Value = "$/saco 25 kilos"
$ws.Range("O53").Value = "Región del Maule"
$ws.Range("P53").Value = 240
$ws.Range("Q53").Value = 25
$ws.Range("R53").Value = "Hortaliza"
